$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1086.4445
$ws.Range("I6").Value = 928
$ws.Range("J6").Value = 1641
$ws.Range("K6").Value = 2784
$ws.Range("L6").Value = 4923
$ws.Range("M6").Value = -2672
$ws.Range("N6").Value = -5147

$ws.Range("H8").Value = 123.875
$ws.Range("I8").Value = 123.875
$ws.Range("K8").Value = 371.625
$ws.Range("M8").Value = -232.625

$ws.Range("H31").Value = 163.25
$ws.Range("I31").Value = 163.25
$ws.Range("K31").Value = 489.75
$ws.Range("M31").Value = -259.75

$ws.Range("H38").Value = 776
$ws.Range("I38").Value = 776
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 2328
$ws.Range("L38").Value = 0
$ws.Range("N38").Value = -1956

$ws.Range("H39").Value = 300
$ws.Range("I39").Value = 300
$ws.Range("K39").Value = 900
$ws.Range("M39").Value = -604

$ws.Range("H88").Value = 3832.6667
$ws.Range("J88").Value = 749
$ws.Range("L88").Value = 749
$ws.Range("N88").Value = -1561

$ws.Range("H91").Value = 3832.6667
$ws.Range("J91").Value = 749
$ws.Range("L91").Value = 749
$ws.Range("N91").Value = -3557

$ws.Range("H92").Value = 811.6316
$ws.Range("I92").Value = 893.53845
$ws.Range("J92").Value = 634.1667
$ws.Range("K92").Value = 893.53845
$ws.Range("L92").Value = 634.1667
$ws.Range("M92").Value = 354.46155
$ws.Range("N92").Value = -3130.1667

$ws.Range("H111").Value = 2023.4445
$ws.Range("I111").Value = 954.1429000000001
$ws.Range("K111").Value = 2862.4287
$ws.Range("M111").Value = 204.5712999999996

$ws.Range("H137").Value = 3160.5264
$ws.Range("J137").Value = 5250
$ws.Range("L137").Value = 15750
$ws.Range("N137").Value = -20850

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7457.3105
$ws.Range("I32").Value = 7457.3105
$ws.Range("K32").Value = 7457.3105
$ws.Range("M32").Value = -7170.3105

$ws.Range("H45").Value = 4580.8
$ws.Range("I45").Value = 4997.5
$ws.Range("J45").Value = 2914
$ws.Range("K45").Value = 4997.5
$ws.Range("L45").Value = 2914
$ws.Range("M45").Value = -4620.5
$ws.Range("N45").Value = -3668

$ws.Range("H102").Value = 3216.6
$ws.Range("I102").Value = 3018.4443
$ws.Range("J102").Value = 5000
$ws.Range("K102").Value = 3018.4443
$ws.Range("L102").Value = 5000
$ws.Range("M102").Value = -1396.4443
$ws.Range("N102").Value = -8244

$ws.Range("H122").Value = 2720.7
$ws.Range("I122").Value = 2800.7778
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 8402.3334
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -5952.3334
$ws.Range("N122").Value = -10900

$ws.Range("H132").Value = 1958.1154
$ws.Range("I132").Value = 1427.909
$ws.Range("K132").Value = 4283.727000000001
$ws.Range("M132").Value = -1753.727000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5626.5713
$ws.Range("I86").Value = 2671.75
$ws.Range("J86").Value = 9566.333000000001
$ws.Range("K86").Value = 2671.75
$ws.Range("L86").Value = 9566.333000000001
$ws.Range("M86").Value = -1548.75
$ws.Range("N86").Value = -11812.333

$ws.Range("H89").Value = 5626.5713
$ws.Range("I89").Value = 2671.75
$ws.Range("J89").Value = 9566.333000000001
$ws.Range("K89").Value = 13358.75
$ws.Range("L89").Value = 47831.665
$ws.Range("M89").Value = -7742.75
$ws.Range("N89").Value = -59063.665

$ws.Range("H94").Value = 3356.9
$ws.Range("I94").Value = 2595
$ws.Range("J94").Value = 4118.8
$ws.Range("K94").Value = 2595
$ws.Range("L94").Value = 4118.8
$ws.Range("M94").Value = -2144
$ws.Range("N94").Value = -5020.8

$ws.Range("H99").Value = 1814.8889
$ws.Range("I99").Value = 1205.6666
$ws.Range("J99").Value = 3033.3333
$ws.Range("K99").Value = 1205.6666
$ws.Range("L99").Value = 3033.3333
$ws.Range("M99").Value = 292.3334
$ws.Range("N99").Value = -6029.3333

$ws.Range("H134").Value = 1735.1538
$ws.Range("I134").Value = 1546.4166
$ws.Range("J134").Value = 4000
$ws.Range("K134").Value = 4639.2498
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = -2104.2498
$ws.Range("N134").Value = -17070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 533.6667
$ws.Range("I16").Value = 294
$ws.Range("K16").Value = 294
$ws.Range("M16").Value = -7

$ws.Range("H113").Value = 533.6667
$ws.Range("I113").Value = 294
$ws.Range("K113").Value = 294
$ws.Range("M113").Value = 1876

$ws.Range("H122").Value = 3643.818
$ws.Range("I122").Value = 2808.2
$ws.Range("J122").Value = 12000
$ws.Range("K122").Value = 8424.599999999999
$ws.Range("L122").Value = 36000
$ws.Range("M122").Value = -5974.599999999999
$ws.Range("N122").Value = -40900

$ws.Range("H132").Value = 1286.909
$ws.Range("I132").Value = 915.7
$ws.Range("K132").Value = 2747.1
$ws.Range("M132").Value = -217.1000000000004

$ws.Range("H134").Value = 2872.95
$ws.Range("I134").Value = 2497.842
$ws.Range("J134").Value = 10000
$ws.Range("K134").Value = 7493.526
$ws.Range("L134").Value = 30000
$ws.Range("M134").Value = -4958.526
$ws.Range("N134").Value = -35070

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 1612.5
$ws.Range("I13").Value = 2550
$ws.Range("J13").Value = 1143.75
$ws.Range("K13").Value = 7650
$ws.Range("L13").Value = 3431.25
$ws.Range("M13").Value = -7482
$ws.Range("N13").Value = -3767.25

$ws.Range("H107").Value = 863.6667
$ws.Range("J107").Value = 863.6667
$ws.Range("L107").Value = 2591.0001
$ws.Range("N107").Value = -6431.0001

$ws.Range("H128").Value = 199997.5
$ws.Range("I128").Value = 199997.5
$ws.Range("K128").Value = 599992.5
$ws.Range("M128").Value = -595012.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 16500
$ws.Range("I80").Value = 3000
$ws.Range("K80").Value = 3000
$ws.Range("M80").Value = -2002

$ws.Range("H83").Value = 16500
$ws.Range("I83").Value = 3000
$ws.Range("K83").Value = 15000
$ws.Range("M83").Value = -10008

$ws.Range("H97").Value = 6823
$ws.Range("I97").Value = 6823
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 6823
$ws.Range("L97").Value = 0
$ws.Range("N97").Value = -6327

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3439.9
$ws.Range("I46").Value = 1800
$ws.Range("J46").Value = 9999.5
$ws.Range("K46").Value = 1800
$ws.Range("L46").Value = 9999.5
$ws.Range("M46").Value = -1612
$ws.Range("N46").Value = -10375.5

$ws.Range("H93").Value = 3459.1667
$ws.Range("I93").Value = 3751.5
$ws.Range("J93").Value = 3313
$ws.Range("K93").Value = 3751.5
$ws.Range("L93").Value = 3313
$ws.Range("M93").Value = -2503.5
$ws.Range("N93").Value = -5809

$ws.Range("H122").Value = 3573.1
$ws.Range("I122").Value = 3878.875
$ws.Range("J122").Value = 2350
$ws.Range("K122").Value = 11636.625
$ws.Range("L122").Value = 7050
$ws.Range("M122").Value = -9186.625
$ws.Range("N122").Value = -11950

$ws.Range("H136").Value = 4999.5
$ws.Range("I136").Value = 4999.5
$ws.Range("K136").Value = 14998.5
$ws.Range("M136").Value = -12448.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 16216.25
$ws.Range("I2").Value = 17104.285
$ws.Range("J2").Value = 10000
$ws.Range("K2").Value = 17104.285
$ws.Range("L2").Value = 10000
$ws.Range("M2").Value = -16992.285
$ws.Range("N2").Value = -10224

$ws.Range("H81").Value = 1849.5
$ws.Range("I81").Value = 1849.5
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 3699
$ws.Range("L81").Value = 0
$ws.Range("N81").Value = -2638

$ws.Range("H84").Value = 1849.5
$ws.Range("I84").Value = 1849.5
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 18495
$ws.Range("L84").Value = 0
$ws.Range("N84").Value = -13191

$ws.Range("H100").Value = 1490.375
$ws.Range("I100").Value = 1022.36365
$ws.Range("J100").Value = 2520
$ws.Range("K100").Value = 2044.7273
$ws.Range("L100").Value = 5040
$ws.Range("M100").Value = -1503.7273
$ws.Range("N100").Value = -6122

$ws.Range("H132").Value = 1631.6666
$ws.Range("J132").Value = 3666
$ws.Range("L132").Value = 10998
$ws.Range("N132").Value = -16058

$ws.Range("H136").Value = 2200.0344
$ws.Range("I136").Value = 2011.7693
$ws.Range("K136").Value = 6035.3079
$ws.Range("M136").Value = -3485.3079
